$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 28
$ws.Cells.Item(3, 6).Value = 50
$ws.Cells.Item(5, 6).Value = 1237
$ws.Cells.Item(6, 6).Value = 1685
$ws.Cells.Item(8, 6).Value = 559
$ws.Cells.Item(9, 6).Value = 2425
$ws.Cells.Item(10, 6).Value = 699
$ws.Cells.Item(11, 6).Value = 562
$ws.Cells.Item(12, 6).Value = 561
$ws.Cells.Item(13, 6).Value = 6
$ws.Cells.Item(14, 6).Value = 643
$ws.Cells.Item(15, 6).Value = 323
$ws.Cells.Item(16, 6).Value = 198
$ws.Cells.Item(19, 6).Value = 1224
$ws.Cells.Item(20, 6).Value = 693
$ws.Cells.Item(22, 6).Value = 2601
$ws.Cells.Item(24, 6).Value = 21
$ws.Cells.Item(28, 6).Value = 1750
$ws.Cells.Item(29, 6).Value = 11
$ws.Cells.Item(31, 6).Value = 516
$ws.Cells.Item(34, 6).Value = 4533
$ws.Cells.Item(35, 6).Value = 94

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 4193
$ws.Cells.Item(8, 6).Value = 58
$ws.Cells.Item(11, 6).Value = 62
$ws.Cells.Item(14, 6).Value = 304
$ws.Cells.Item(20, 6).Value = 270
$ws.Cells.Item(26, 6).Value = 228
$ws.Cells.Item(27, 6).Value = 8
$ws.Cells.Item(28, 6).Value = 247

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 500
$ws.Cells.Item(7, 6).Value = 165

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 500
$ws.Cells.Item(7, 6).Value = 28
$ws.Cells.Item(8, 6).Value = 50
$ws.Cells.Item(10, 6).Value = 1237
$ws.Cells.Item(11, 6).Value = 1685
$ws.Cells.Item(13, 6).Value = 58
$ws.Cells.Item(16, 6).Value = 559
$ws.Cells.Item(17, 6).Value = 2425
$ws.Cells.Item(18, 6).Value = 699
$ws.Cells.Item(19, 6).Value = 562
$ws.Cells.Item(20, 6).Value = 561
$ws.Cells.Item(21, 6).Value = 643
$ws.Cells.Item(22, 6).Value = 323
$ws.Cells.Item(23, 6).Value = 62
$ws.Cells.Item(24, 6).Value = 198
$ws.Cells.Item(28, 6).Value = 1224
$ws.Cells.Item(29, 6).Value = 693
$ws.Cells.Item(32, 6).Value = 2601
$ws.Cells.Item(33, 6).Value = 270
$ws.Cells.Item(35, 6).Value = 21
$ws.Cells.Item(38, 6).Value = 165
$ws.Cells.Item(41, 6).Value = 1750
$ws.Cells.Item(42, 6).Value = 228
$ws.Cells.Item(44, 6).Value = 516
$ws.Cells.Item(47, 6).Value = 4533
$ws.Cells.Item(48, 6).Value = 94
